# Fix material database field line traces example:
# correct the header labels in row 1 of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "material name"
$ws.Range("B1").Value = "Conductivity"
$ws.Range("C1").Value = "permittivity"

$ws.Range("C1").Select()
